# Add CAN bus transceiver (TCAN1042HVDRQ1) line item to the BOM.
#
# Before: row 28 is the trailing blank placeholder row (with G28 holding the
# shared SUM-style formula "=F28*E28" but no data).
# After:  row 28 is filled in with the new part, and a new blank placeholder
# row 29 is appended below it (inheriting row 28's old formatting/formula).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Create the new trailing blank row 29 by copying row 28's current
#    formatting (borders / wrap / hyperlink style / etc.) down a row, then
#    give it the same "price * qty" formula the old blank row had.
$ws.Range("A28:G28").Copy()
$ws.Range("A29:G29").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G29").Formula = "=F29*E29"

# 2) Populate row 28 with the new CAN transceiver part.
#    (Manufacturer part number first, then description, matches the order
#    the new shared-string entries were originally authored in.)
$ws.Range("C28").Value = "TCAN1042HVDRQ1"
$ws.Range("B28").Value = "IC TXRX CAN FAULT PROT 8SOIC"
$ws.Hyperlinks.Add($ws.Range("D28"), "http://www.ti.com/lit/ds/symlink/tcan1042hgv-q1.pdf") | Out-Null
$ws.Range("E28").Value = 2.05
$ws.Range("F28").Value = 1
$ws.Range("G28").Formula = "=F28*E28"

# Adding the hyperlink via Hyperlinks.Add() re-styles D28 with Excel's
# built-in "Hyperlink" cell style; restore the sheet's existing
# bordered-hyperlink look (copied from the still-intact D27 cell) instead.
$ws.Range("D27").Copy()
$ws.Range("D28").PasteSpecial(-4122)        # xlPasteFormats

# 3) Match the saved selection state (active cell on the new blank row).
$ws.Range("E29").Select() | Out-Null
